$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 48.177
$ws.Range("D2").Value = 48.177
$ws.Range("E2").Value = 3.55277867
$ws.Range("F2").Value = 0.02433503
$ws.Range("G2").Value = 1.164731
$ws.Range("H2").Value = 56.44339078
$ws.Range("I2").Value = 5.566206396968554
$ws.Range("J2").Value = 5.566206396968554
$ws.Range("K2").Value = 0.395659608027096
$ws.Range("L2").Value = 0.003768716893868008
$ws.Range("M2").Value = 0.1840152679911476
$ws.Range("N2").Value = 12.67829446369795
$ws.Range("C3").Value = 86.72
$ws.Range("D3").Value = 86.72
$ws.Range("E3").Value = 1.99003187
$ws.Range("F3").Value = 0.01624326
$ws.Range("G3").Value = 1.38270517
$ws.Range("H3").Value = 120.22189381
$ws.Range("I3").Value = 12.77132360825246
$ws.Range("J3").Value = 12.77132360825246
$ws.Range("K3").Value = 0.287110909693568
$ws.Range("L3").Value = 0.002566627816951798
$ws.Range("M3").Value = 0.1272442274548896
$ws.Range("N3").Value = 22.89348009974907
$ws.Range("C4").Value = 28.107
$ws.Range("D4").Value = 56.154
$ws.Range("E4").Value = 3.11514325
$ws.Range("F4").Value = 0.03845762
$ws.Range("G4").Value = 0.53157267
$ws.Range("H4").Value = 15.22160484
$ws.Range("I4").Value = 5.464193298960828
$ws.Range("J4").Value = 10.8941712381882
$ws.Range("K4").Value = 0.5607040403264237
$ws.Range("L4").Value = 0.006600256731713133
$ws.Range("M4").Value = 0.1009696213690559
$ws.Range("N4").Value = 5.29275551823316
$ws.Range("C5").Value = 48.201
$ws.Range("D5").Value = 94.00700000000001
$ws.Range("E5").Value = 1.84563495
$ws.Range("F5").Value = 0.02635705
$ws.Range("G5").Value = 0.62059492
$ws.Range("H5").Value = 30.10493904
$ws.Range("I5").Value = 8.487086050561148
$ws.Range("J5").Value = 15.4521752747217
$ws.Range("K5").Value = 0.2984568981797348
$ws.Range("L5").Value = 0.005955568334902362
$ws.Range("M5").Value = 0.1131074828464919
$ws.Range("N5").Value = 8.3538971929861
$ws.Range("C6").Value = 16.658
$ws.Range("D6").Value = 66.502
$ws.Range("E6").Value = 2.65270173
$ws.Range("F6").Value = 0.05174844
$ws.Range("G6").Value = 0.2110295
$ws.Range("H6").Value = 3.59561762
$ws.Range("I6").Value = 3.661207511250959
$ws.Range("J6").Value = 14.57815209546527
$ws.Range("K6").Value = 0.532946189282725
$ws.Range("L6").Value = 0.01399910313752474
$ws.Range("M6").Value = 0.05915484469594186
$ws.Range("N6").Value = 1.539338216296044
$ws.Range("C7").Value = 26.79
$ws.Range("D7").Value = 98.184
$ws.Range("E7").Value = 1.76667751
$ws.Range("F7").Value = 0.03642079
$ws.Range("G7").Value = 0.23882021
$ws.Range("H7").Value = 6.50948818
$ws.Range("I7").Value = 5.326750817738171
$ws.Range("J7").Value = 16.06113645120121
$ws.Range("K7").Value = 0.285562182767128
$ws.Range("L7").Value = 0.007647379164392464
$ws.Range("M7").Value = 0.04746934300979614
$ws.Range("N7").Value = 2.415783426311956
$ws.Range("C8").Value = 11.234
$ws.Range("D8").Value = 67.215
$ws.Range("E8").Value = 2.6200787
$ws.Range("F8").Value = 0.0551564
$ws.Range("G8").Value = 0.10283655
$ws.Range("H8").Value = 1.19411918
$ws.Range("I8").Value = 2.2504018104126
$ws.Range("J8").Value = 13.45086245192853
$ws.Range("K8").Value = 0.553914384164642
$ws.Range("L8").Value = 0.01429440084562958
$ws.Range("M8").Value = 0.03249865855701594
$ws.Range("N8").Value = 0.5401764753912315
$ws.Range("C9").Value = 19.057
$ws.Range("D9").Value = 95.72
$ws.Range("E9").Value = 1.81713302
$ws.Range("F9").Value = 0.04022206
$ws.Range("G9").Value = 0.12506851
$ws.Range("H9").Value = 2.45405086
$ws.Range("I9").Value = 4.416475546887614
$ws.Range("J9").Value = 16.03292408310644
$ws.Range("K9").Value = 0.3185761461047367
$ws.Range("L9").Value = 0.009349877505388279
$ws.Range("M9").Value = 0.03147370869819951
$ws.Range("N9").Value = 1.095591522548135
$ws.Range("C10").Value = 8.576000000000001
$ws.Range("D10").Value = 68.38
$ws.Range("E10").Value = 2.62732248
$ws.Range("F10").Value = 0.05605742
$ws.Range("G10").Value = 0.06064484000000001
$ws.Range("H10").Value = 0.55398124
$ws.Range("I10").Value = 2.054859227890428
$ws.Range("J10").Value = 16.34842105130862
$ws.Range("K10").Value = 0.6844131339053274
$ws.Range("L10").Value = 0.01744708754512934
$ws.Range("M10").Value = 0.0260100294426209
$ws.Range("N10").Value = 0.3423370118917518
$ws.Range("C11").Value = 14.599
$ws.Range("D11").Value = 89.194
$ws.Range("E11").Value = 1.96091773
$ws.Range("F11").Value = 0.04334522
$ws.Range("G11").Value = 0.07787151
$ws.Range("H11").Value = 1.19065976
$ws.Range("I11").Value = 3.72075839622019
$ws.Range("J11").Value = 16.57845617689443
$ws.Range("K11").Value = 0.3671044303728115
$ws.Range("L11").Value = 0.01189198948727948
$ws.Range("M11").Value = 0.02623187871428161
$ws.Range("N11").Value = 0.6772005714042574
$ws.Range("C12").Value = 6.887
$ws.Range("D12").Value = 68.515
$ws.Range("E12").Value = 2.6294566
$ws.Range("F12").Value = 0.05679174000000001
$ws.Range("G12").Value = 0.04003063
$ws.Range("H12").Value = 0.29944527
$ws.Range("I12").Value = 1.72545883990555
$ws.Range("J12").Value = 17.11685607118789
$ws.Range("K12").Value = 0.6874198937515303
$ws.Range("L12").Value = 0.0204690165765728
$ws.Range("M12").Value = 0.02060578349436188
$ws.Range("N12").Value = 0.2249043335339928
$ws.Range("C13").Value = 11.917
$ws.Range("D13").Value = 81.919
$ws.Range("E13").Value = 2.14064073
$ws.Range("F13").Value = 0.03880975
$ws.Range("G13").Value = 0.04558601
$ws.Range("H13").Value = 0.5788314899999999
$ws.Range("I13").Value = 3.396416736304107
$ws.Range("J13").Value = 15.96663309443705
$ws.Range("K13").Value = 0.4127071059568448
$ws.Range("L13").Value = 0.01067296300831577
$ws.Range("M13").Value = 0.0167084814250557
$ws.Range("N13").Value = 0.4064593869378709
